$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2116040955631399
$ws.Range("C2").Value = 0.5290102389078498
$ws.Range("J2").Value = 0.0136518771331058
$ws.Range("P2").Value = 0.1399317406143345
$ws.Range("S2").Value = 0.10580204778157
$ws.Range("C3").Value = 0.01257861635220126
$ws.Range("J3").Value = 0.03773584905660377
$ws.Range("P3").Value = 0.6918238993710691
$ws.Range("S3").Value = 0.2578616352201258
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.6585365853658537
$ws.Range("S4").Value = 0.2682926829268293
$ws.Range("B6").Value = 0.06046511627906977
$ws.Range("F6").Value = 0.05116279069767442
$ws.Range("J6").Value = 0.2930232558139535
$ws.Range("O6").Value = 0.03720930232558139
$ws.Range("Q6").Value = 0.1906976744186047
$ws.Range("R6").Value = 0.05116279069767442
$ws.Range("S6").Value = 0.3162790697674419
$ws.Range("B7").Value = 0.11328125
$ws.Range("D7").Value = 0.015625
$ws.Range("E7").Value = 0.00390625
$ws.Range("F7").Value = 0.05859375
$ws.Range("J7").Value = 0.09765625
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.2109375
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.40625
$ws.Range("B8").Value = 0.07575757575757576
$ws.Range("D8").Value = 0.01515151515151515
$ws.Range("E8").Value = 0.002164502164502165
$ws.Range("F8").Value = 0.05627705627705628
$ws.Range("J8").Value = 0.119047619047619
$ws.Range("O8").Value = 0.02597402597402598
$ws.Range("Q8").Value = 0.2424242424242424
$ws.Range("R8").Value = 0.07792207792207792
$ws.Range("S8").Value = 0.3852813852813853
$ws.Range("B9").Value = 0.07253886010362694
$ws.Range("D9").Value = 0.02072538860103627
$ws.Range("F9").Value = 0.04145077720207254
$ws.Range("J9").Value = 0.1036269430051813
$ws.Range("O9").Value = 0.005181347150259068
$ws.Range("Q9").Value = 0.2227979274611399
$ws.Range("R9").Value = 0.08808290155440414
$ws.Range("S9").Value = 0.4455958549222798
$ws.Range("B10").Value = 0.09970887918486172
$ws.Range("D10").Value = 0.02037845705967977
$ws.Range("E10").Value = 0.002911208151382824
$ws.Range("F10").Value = 0.05604075691411936
$ws.Range("J10").Value = 0.1491994177583697
$ws.Range("O10").Value = 0.01746724890829694
$ws.Range("Q10").Value = 0.2241630276564774
$ws.Range("R10").Value = 0.07132459970887918
$ws.Range("S10").Value = 0.3588064046579331
$ws.Range("G11").Value = 0.154639175257732
$ws.Range("J11").Value = 0.08762886597938144
$ws.Range("K11").Value = 0.2164948453608248
$ws.Range("L11").Value = 0.5309278350515464
$ws.Range("S11").Value = 0.01030927835051546
$ws.Range("G12").Value = 0.7819905213270142
$ws.Range("J12").Value = 0.1753554502369668
$ws.Range("L12").Value = 0.02369668246445497
$ws.Range("S12").Value = 0.01895734597156398
$ws.Range("G13").Value = 0.7333333333333333
$ws.Range("J13").Value = 0.1833333333333333
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.00796812749003984
$ws.Range("H15").Value = 0.1752988047808765
$ws.Range("I15").Value = 0.05179282868525897
$ws.Range("J15").Value = 0.3386454183266932
$ws.Range("K15").Value = 0.0398406374501992
$ws.Range("M15").Value = 0.0199203187250996
$ws.Range("O15").Value = 0.08764940239043825
$ws.Range("S15").Value = 0.2788844621513944
$ws.Range("F16").Value = 0.02312138728323699
$ws.Range("H16").Value = 0.1849710982658959
$ws.Range("I16").Value = 0.07514450867052024
$ws.Range("J16").Value = 0.3872832369942196
$ws.Range("K16").Value = 0.1213872832369942
$ws.Range("M16").Value = 0.03468208092485549
$ws.Range("O16").Value = 0.04046242774566474
$ws.Range("S16").Value = 0.1329479768786127
$ws.Range("F17").Value = 0.01615798922800718
$ws.Range("H17").Value = 0.1651705565529623
$ws.Range("I17").Value = 0.0843806104129264
$ws.Range("J17").Value = 0.3680430879712747
$ws.Range("K17").Value = 0.1274685816876122
$ws.Range("M17").Value = 0.03590664272890485
$ws.Range("O17").Value = 0.0718132854578097
$ws.Range("S17").Value = 0.1310592459605027
$ws.Range("F18").Value = 0.03910614525139665
$ws.Range("H18").Value = 0.1675977653631285
$ws.Range("I18").Value = 0.1173184357541899
$ws.Range("J18").Value = 0.3798882681564246
$ws.Range("K18").Value = 0.1173184357541899
$ws.Range("M18").Value = 0.01675977653631285
$ws.Range("O18").Value = 0.0670391061452514
$ws.Range("S18").Value = 0.09497206703910614
$ws.Range("F19").Value = 0.02075611564121571
$ws.Range("H19").Value = 0.195700518902891
$ws.Range("I19").Value = 0.07338769458858414
$ws.Range("J19").Value = 0.3750926612305411
$ws.Range("K19").Value = 0.1297257227575982
$ws.Range("M19").Value = 0.02001482579688658
$ws.Range("O19").Value = 0.06597479614529281
$ws.Range("S19").Value = 0.1193476649369904

Write-Output "Applied team specific time data updates"
